$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.532.63"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.472.29"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("E5").Value = "  -3.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3062"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.055"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06620"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.467"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001029"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "1.473.60"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05895"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9634"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.466"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.247"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").Value = "20.587.72"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.128"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "1.631.38"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.978"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8135"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.966"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07942"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.535"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.223"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05796"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.710"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9590"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.594"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1878"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5269"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5175"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.789"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06459"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9944"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
